$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update the first four sample-size values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON): D2 value was deleted, B2 and C2 now carry values instead
$ws.Range("D2").Value = ""
$ws.Range("B2").Value = 27.667255443025404
$ws.Range("C2").Value = 26.086840642673973

# Row 3 (STR): B3 value was deleted, C3 value updated
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = 24.120046273807962

# Selection now only covers B1:E3 instead of the full B1:AY3 range
$ws.Range("B1:E3").Select()
